# Repull data, push all data, mean calculation
# Update column F (dSF) values on several rows to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -3
    6  = -3
    10 = 3
    11 = -1
    14 = 3
    16 = 1
    26 = 0
    32 = -1
    34 = 2
    41 = 3
    47 = 8
    48 = 0
    49 = 5
    51 = 1
    55 = -2
    58 = 3
    60 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
